# Refresh "need_to_buy" forecast data (rows 2-15) with the latest values from R.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 1).Value = 45936
$ws.Cells.Item(2, 2).Value = 4444.43724989551
$ws.Cells.Item(2, 3).Value = 5577.9838050625
$ws.Cells.Item(2, 4).Value = 3620
$ws.Cells.Item(2, 5).Value = 6433.405962
$ws.Cells.Item(2, 6).Value = 164.456354881958

# Row 3
$ws.Cells.Item(3, 1).Value = 45937
$ws.Cells.Item(3, 2).Value = 4428.72629889378
$ws.Cells.Item(3, 3).Value = 5532.40196082231
$ws.Cells.Item(3, 4).Value = 3620
$ws.Cells.Item(3, 5).Value = 6167.167887
$ws.Cells.Item(3, 6).Value = 152.118481205356

# Row 4
$ws.Cells.Item(4, 1).Value = 45938
$ws.Cells.Item(4, 2).Value = 4362.60881602614
$ws.Cells.Item(4, 3).Value = 5352.72560003091
$ws.Cells.Item(4, 4).Value = 3620
$ws.Cells.Item(4, 5).Value = 6041.754981
$ws.Cells.Item(4, 6).Value = 142.161323541865

# Row 5
$ws.Cells.Item(5, 1).Value = 45939
$ws.Cells.Item(5, 2).Value = 4321.34845690586
$ws.Cells.Item(5, 3).Value = 5169.61121455699
$ws.Cells.Item(5, 4).Value = 3620
$ws.Cells.Item(5, 5).Value = 5961.829301
$ws.Cells.Item(5, 6).Value = 132.920502443797

# Row 6
$ws.Cells.Item(6, 1).Value = 45940
$ws.Cells.Item(6, 2).Value = 5297.0753414716
$ws.Cells.Item(6, 3).Value = 4688.84249212959
$ws.Cells.Item(6, 4).Value = 3620
$ws.Cells.Item(6, 5).Value = 7641.295045
$ws.Cells.Item(6, 6).Value = 142.210924819083

# Row 7
$ws.Cells.Item(7, 1).Value = 45941
$ws.Cells.Item(7, 2).Value = 1503.60506415849
$ws.Cells.Item(7, 3).Value = 2649.53609266902
$ws.Cells.Item(7, 4).Value = 3620
$ws.Cells.Item(7, 5).Value = 3877.880684
$ws.Cells.Item(7, 6).Value = 58.4921546879391

# Row 8
$ws.Cells.Item(8, 1).Value = 45942
$ws.Cells.Item(8, 2).Value = 1399.99543622148
$ws.Cells.Item(8, 3).Value = 2751.86219463129
$ws.Cells.Item(8, 4).Value = 3620
$ws.Cells.Item(8, 5).Value = 3764.615186
$ws.Cells.Item(8, 6).Value = 62.3534143504088

# Row 9
$ws.Cells.Item(9, 1).Value = 45943
$ws.Cells.Item(9, 2).Value = 5841.28374289435
$ws.Cells.Item(9, 3).Value = 5608.97789339652
$ws.Cells.Item(9, 4).Value = 3620
$ws.Cells.Item(9, 5).Value = 8951.246298
$ws.Cells.Item(9, 6).Value = 212.455852020924

# Row 10
$ws.Cells.Item(10, 1).Value = 45944
$ws.Cells.Item(10, 2).Value = 5841.28374289435
$ws.Cells.Item(10, 3).Value = 5734.16455495641
$ws.Cells.Item(10, 4).Value = 3620
$ws.Cells.Item(10, 5).Value = 8951.246298
$ws.Cells.Item(10, 6).Value = 217.671962919252

# Row 11
$ws.Cells.Item(11, 1).Value = 45945
$ws.Cells.Item(11, 2).Value = 5841.28374289435
$ws.Cells.Item(11, 3).Value = 5822.82193336872
$ws.Cells.Item(11, 4).Value = 3620
$ws.Cells.Item(11, 5).Value = 8970.797255
$ws.Cells.Item(11, 6).Value = 222.180643561432

# Row 12
$ws.Cells.Item(12, 1).Value = 45946
$ws.Cells.Item(12, 2).Value = 5841.28374289435
$ws.Cells.Item(12, 3).Value = 6002.33643971666
$ws.Cells.Item(12, 4).Value = 3620
$ws.Cells.Item(12, 5).Value = 8970.797255
$ws.Cells.Item(12, 6).Value = 229.660414659263

# Row 13
$ws.Cells.Item(13, 1).Value = 45947
$ws.Cells.Item(13, 2).Value = 5841.28374289435
$ws.Cells.Item(13, 3).Value = 5443.05670487601
$ws.Cells.Item(13, 4).Value = 3620
$ws.Cells.Item(13, 5).Value = 8970.797255
$ws.Cells.Item(13, 6).Value = 206.357092374236

# Row 14
$ws.Cells.Item(14, 1).Value = 45948
$ws.Cells.Item(14, 2).Value = 1742.27770790123
$ws.Cells.Item(14, 3).Value = 3422.62263316978
$ws.Cells.Item(14, 4).Value = 3620
$ws.Cells.Item(14, 5).Value = 4473.64226
$ws.Cells.Item(14, 6).Value = 105.582799386189

# Row 15
$ws.Cells.Item(15, 1).Value = 45949
$ws.Cells.Item(15, 2).Value = 1636.94065696827
$ws.Cells.Item(15, 3).Value = 3526.01895954624
$ws.Cells.Item(15, 4).Value = 3620
$ws.Cells.Item(15, 5).Value = 4358.481159
$ws.Cells.Item(15, 6).Value = 109.481644232415

